# Apply cryptos list update per commit "Updated cryptos list on Wed Feb 14 06:27:48 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Range, [string]$Text)
    # Force literal text (avoids Excel auto-converting numeric-looking
    # strings like "69.00" into the number 69) while preserving the
    # cell's original style/number format.
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

$ws.Range("D2").Value = "49.610.92"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.645.91"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCellValue $ws.Range("D5") "112.64"
$ws.Range("E5").Value = "  -1.29%  "
Set-TextCellValue $ws.Range("D6") "326.19"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -3.16%  "
Set-TextCellValue $ws.Range("D11") "20.01"
$ws.Range("E11").Value = "  -0.92%  "
Set-TextCellValue $ws.Range("D12") "0.0814"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  +1.90%  "
Set-TextCellValue $ws.Range("D14") "7.56"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").Value = "3.057.64"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "2.641.95"
$ws.Range("E16").Value = "  -0.81%  "
Set-TextCellValue $ws.Range("D17") "0.860"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "49.574.41"
$ws.Range("E18").Value = "  -0.73%  "
Set-TextCellValue $ws.Range("D19") "13.34"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCellValue $ws.Range("D20") "6.68"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCellValue $ws.Range("D21") "2.90"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -1.03%  "
Set-TextCellValue $ws.Range("D23") "268.31"
$ws.Range("E23").Value = "  -3.19%  "
Set-TextCellValue $ws.Range("D24") "69.00"
$ws.Range("E24").Value = "  -4.30%  "
Set-TextCellValue $ws.Range("D25") "2.57"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCellValue $ws.Range("D26") "26.10"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCellValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  -1.11%  "
Set-TextCellValue $ws.Range("D30") "0.138"
$ws.Range("E30").Value = "  -2.42%  "
Set-TextCellValue $ws.Range("D31") "34.74"
$ws.Range("E31").Value = "  -3.76%  "
Set-TextCellValue $ws.Range("D32") "49.65"
$ws.Range("E32").Value = "  -1.25%  "
Set-TextCellValue $ws.Range("D33") "5.47"
$ws.Range("E33").Value = "  +0.25%  "
Set-TextCellValue $ws.Range("D34") "0.0822"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCellValue $ws.Range("D35") "19.17"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCellValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -2.03%  "
Set-TextCellValue $ws.Range("D39") "3.12"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCellValue $ws.Range("D40") "129.61"
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCellValue $ws.Range("D41") "23.68"
$ws.Range("E41").Value = "  +7.82%  "
$ws.Range("E42").Value = "  +2.94%  "
Set-TextCellValue $ws.Range("D43") "0.0343"
$ws.Range("E43").Value = "  +8.95%  "
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "2.061.04"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  -0.84%  "
Set-TextCellValue $ws.Range("D47") "2.12"
$ws.Range("E47").Value = "  +6.49%  "
$ws.Range("E48").Value = "  -4.97%  "
Set-TextCellValue $ws.Range("D49") "8.92"
$ws.Range("E49").Value = "  -2.59%  "
Set-TextCellValue $ws.Range("D50") "5.25"
$ws.Range("E50").Value = "  -2.81%  "
Set-TextCellValue $ws.Range("D51") "58.76"
$ws.Range("E51").Value = "  -1.93%  "
